$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.908.60"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.844.71"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "309.47"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4762"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.63%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3667"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.95%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07202"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9271"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.12%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "19.70"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.83%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07702"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "1.874.82"
$ws.Range("E13").Value = "  +3.44%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.316"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.17%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.406"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "88.79"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("E17").Value = "  -0.04%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008640"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.07%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.57"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "26.929.65"
$ws.Range("E21").Value = "  +1.15%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.053"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.84%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.64"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.920"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.36"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "18.16"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.005"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.78%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "114.26"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.927"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.66%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.08876"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.309"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.40%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.176"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.99%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.7490"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.485"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.17%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.734"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.094"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.69%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01955"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.80%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05262"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("E39").Value = "  +1.78%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5203"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.19%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.962"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.50%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1510"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.224"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.44%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "10.51"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.95%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4729"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.87%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.008"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "101.54"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.46%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.603"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.08%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "66.08"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.84%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06024"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.79%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.8866"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.31%  "
